# Auto-generated edit script applying scheduled market-data refresh updates
# to the per-job "Leve Profits" worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each touched row holds cached market-board figures (columns H-N) for one Leve;
# this mirrors a scheduled runner re-pulling current prices and recalculating profit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 639.44446
$ws.Range("I19").Value = 718.2727
$ws.Range("K19").Value = 718.2727
$ws.Range("M19").Value = -543.2727
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H138").Value = 2215.9333
$ws.Range("I138").Value = 1386.909
$ws.Range("J138").Value = 4495.75
$ws.Range("K138").Value = 4160.727000000001
$ws.Range("L138").Value = 13487.25
$ws.Range("M138").Value = 979.2729999999992
$ws.Range("N138").Value = -23767.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 169.28572
$ws.Range("I4").Value = 189.16667
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 189.16667
$ws.Range("L4").Value = 50
$ws.Range("M4").Value = -73.16667000000001
$ws.Range("N4").Value = -282
$ws.Range("H8").Value = 3750
$ws.Range("J8").Value = 7000
$ws.Range("L8").Value = 7000
$ws.Range("N8").Value = -7288
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H13").Value = 2998.5
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 2998.5
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 2998.5
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -3286.5
$ws.Range("H32").Value = 13010.852
$ws.Range("I32").Value = 7173.3687
$ws.Range("K32").Value = 7173.3687
$ws.Range("M32").Value = -6886.3687
$ws.Range("H63").Value = 5302.5713
$ws.Range("I63").Value = 3039.6
$ws.Range("K63").Value = 3039.6
$ws.Range("M63").Value = -2353.6
$ws.Range("H66").Value = 5302.5713
$ws.Range("I66").Value = 3039.6
$ws.Range("K66").Value = 15198
$ws.Range("M66").Value = -11766
$ws.Range("H74").Value = 2529.44
$ws.Range("I74").Value = 2392.762
$ws.Range("K74").Value = 2392.762
$ws.Range("M74").Value = -1518.762
$ws.Range("H77").Value = 2529.44
$ws.Range("I77").Value = 2392.762
$ws.Range("K77").Value = 11963.81
$ws.Range("M77").Value = -7595.810000000001
$ws.Range("H88").Value = 3285.7144
$ws.Range("I88").Value = 4566.6665
$ws.Range("J88").Value = 2325
$ws.Range("K88").Value = 4566.6665
$ws.Range("L88").Value = 2325
$ws.Range("M88").Value = -4160.6665
$ws.Range("N88").Value = -3137
$ws.Range("H91").Value = 3285.7144
$ws.Range("I91").Value = 4566.6665
$ws.Range("J91").Value = 2325
$ws.Range("K91").Value = 4566.6665
$ws.Range("L91").Value = 2325
$ws.Range("M91").Value = -3162.6665
$ws.Range("N91").Value = -5133
$ws.Range("H96").Value = 15736
$ws.Range("J96").Value = 15736
$ws.Range("L96").Value = 15736
$ws.Range("N96").Value = -21228
$ws.Range("H97").Value = 1604.125
$ws.Range("I97").Value = 1805.5
$ws.Range("K97").Value = 1805.5
$ws.Range("M97").Value = -1309.5
$ws.Range("H102").Value = 4381.524
$ws.Range("I102").Value = 2134.4
$ws.Range("J102").Value = 9999.333000000001
$ws.Range("K102").Value = 2134.4
$ws.Range("L102").Value = 9999.333000000001
$ws.Range("M102").Value = -512.4000000000001
$ws.Range("N102").Value = -13243.333
$ws.Range("H137").Value = 69999
$ws.Range("J137").Value = 69999
$ws.Range("L137").Value = 69999
$ws.Range("N137").Value = -80199

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 94999
$ws.Range("J59").Value = 94999
$ws.Range("L59").Value = 94999
$ws.Range("N59").Value = -96693
$ws.Range("H134").Value = 4218
$ws.Range("I134").Value = 1805.8334
$ws.Range("K134").Value = 5417.5002
$ws.Range("M134").Value = -2882.5002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6365.6113
$ws.Range("I31").Value = 4973
$ws.Range("J31").Value = 14999.8
$ws.Range("K31").Value = 4973
$ws.Range("L31").Value = 14999.8
$ws.Range("M31").Value = -4678
$ws.Range("N31").Value = -15589.8
$ws.Range("H34").Value = 6365.6113
$ws.Range("I34").Value = 4973
$ws.Range("J34").Value = 14999.8
$ws.Range("K34").Value = 4973
$ws.Range("L34").Value = 14999.8
$ws.Range("M34").Value = -4771
$ws.Range("N34").Value = -15403.8
$ws.Range("H58").Value = 1920.8823
$ws.Range("I58").Value = 1920.8823
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1920.8823
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1717.8823
$ws.Range("N58").ClearContents()
$ws.Range("H108").Value = 49799.5
$ws.Range("J108").Value = 49799.5
$ws.Range("L108").Value = 49799.5
$ws.Range("N108").Value = -57479.5
$ws.Range("H136").Value = 1920.8823
$ws.Range("I136").Value = 1920.8823
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5762.6469
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3212.6469
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 326.75
$ws.Range("I8").Value = 326.75
$ws.Range("K8").Value = 980.25
$ws.Range("M8").Value = -841.25
$ws.Range("H14").Value = 1997.8182
$ws.Range("I14").Value = 1997.8182
$ws.Range("K14").Value = 5993.4546
$ws.Range("M14").Value = -5820.4546
$ws.Range("H101").Value = 12000
$ws.Range("J101").Value = 12000
$ws.Range("L101").Value = 36000
$ws.Range("N101").Value = -40868

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 16667
$ws.Range("I40").Value = 16667
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 16667
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -16516
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9152.267
$ws.Range("I40").Value = 8523.833000000001
$ws.Range("K40").Value = 8523.833000000001
$ws.Range("M40").Value = -8387.833000000001
$ws.Range("H82").Value = 3729.4
$ws.Range("I82").Value = 3539.75
$ws.Range("J82").Value = 4488
$ws.Range("K82").Value = 3539.75
$ws.Range("L82").Value = 4488
$ws.Range("M82").Value = -3178.75
$ws.Range("N82").Value = -5210
$ws.Range("H85").Value = 3729.4
$ws.Range("I85").Value = 3539.75
$ws.Range("J85").Value = 4488
$ws.Range("K85").Value = 3539.75
$ws.Range("L85").Value = 4488
$ws.Range("M85").Value = -2291.75
$ws.Range("N85").Value = -6984
$ws.Range("H93").Value = 1515.5264
$ws.Range("I93").Value = 1497.4
$ws.Range("K93").Value = 1497.4
$ws.Range("M93").Value = -249.4000000000001
$ws.Range("H100").Value = 6789.0527
$ws.Range("I100").Value = 3570.2856
$ws.Range("J100").Value = 8666.666999999999
$ws.Range("K100").Value = 3570.2856
$ws.Range("L100").Value = 8666.666999999999
$ws.Range("M100").Value = -3029.2856
$ws.Range("N100").Value = -9748.666999999999
$ws.Range("H132").Value = 6448.4287
$ws.Range("I132").Value = 4842.143
$ws.Range("K132").Value = 14526.429
$ws.Range("M132").Value = -11996.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 10293.667
$ws.Range("I4").Value = 13153.125
$ws.Range("J4").Value = 4574.75
$ws.Range("K4").Value = 13153.125
$ws.Range("L4").Value = 4574.75
$ws.Range("M4").Value = -13040.125
$ws.Range("N4").Value = -4800.75
$ws.Range("H96").Value = 2177
$ws.Range("I96").Value = 1295
$ws.Range("J96").Value = 3500
$ws.Range("K96").Value = 1295
$ws.Range("L96").Value = 3500
$ws.Range("M96").Value = 78
$ws.Range("N96").Value = -6246
$ws.Range("H132").Value = 7984.6665
$ws.Range("I132").Value = 7977
$ws.Range("K132").Value = 23931
$ws.Range("M132").Value = -21401
